$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two question cells with the new (longer) text
$ws.Range("A3").Value = "How are you dear?"
$ws.Range("A2").Value = " What is ur name dear?"

# Move the active selection to A2 (matches saved sheetView state)
$ws.Range("A2").Select()
